$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '88.863.15'
$ws.Cells.Item(2, 5).Value = '  +1.18%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.266.66'
$ws.Cells.Item(3, 5).Value = '  -2.53%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.13%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '211.22'
$ws.Cells.Item(5, 5).Value = '  -4.05%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '624.79'
$ws.Cells.Item(6, 5).Value = '  -2.43%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.372'
$ws.Cells.Item(7, 5).Value = '  +15.10%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.714'
$ws.Cells.Item(8, 5).Value = '  +15.62%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.999'
$ws.Cells.Item(9, 5).Value = '  +0.07%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '3.263.72'
$ws.Cells.Item(10, 5).Value = '  -2.79%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.572'
$ws.Cells.Item(11, 5).Value = '  -6.54%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.186'
$ws.Cells.Item(12, 5).Value = '  +10.94%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -5.97%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Toncoin'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.44'
$ws.Cells.Item(14, 5).Value = '  +0.41%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.864.61'
$ws.Cells.Item(15, 5).Value = '  -2.44%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '33.86'
$ws.Cells.Item(16, 5).Value = '  -1.98%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '88.634.59'
$ws.Cells.Item(17, 5).Value = '  +1.37%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.277.75'
$ws.Cells.Item(18, 5).Value = '  -2.01%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.14'
$ws.Cells.Item(19, 5).Value = '  -2.62%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '13.98'
$ws.Cells.Item(20, 5).Value = '  -5.02%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '434.58'
$ws.Cells.Item(21, 5).Value = '  -3.68%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '8.84'
$ws.Cells.Item(22, 5).Value = '  -4.01%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.30'
$ws.Cells.Item(23, 5).Value = '  -0.65%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '7.37'
$ws.Cells.Item(24, 5).Value = '  -0.50%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '5.20'
$ws.Cells.Item(25, 5).Value = '  -3.35%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '12.14'
$ws.Cells.Item(26, 5).Value = '  -1.82%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '3.461.93'
$ws.Cells.Item(27, 5).Value = '  -1.46%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '76.64'
$ws.Cells.Item(28, 5).Value = '  -3.05%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +2.46%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  +0.02%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.179'
$ws.Cells.Item(31, 5).Value = '  -4.97%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.37%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Bittensor'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '558.91'
$ws.Cells.Item(33, 5).Value = '  -7.60%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '8.69'
$ws.Cells.Item(34, 5).Value = '  -7.00%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.37'
$ws.Cells.Item(35, 5).Value = '  -12.55%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -5.29%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '7.07'
$ws.Cells.Item(37, 5).Value = '  +6.41%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -8.66%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '22.58'
$ws.Cells.Item(39, 5).Value = '  -4.20%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '21.85'
$ws.Cells.Item(40, 5).Value = '  +2.12%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 5).Value = '  +0.07%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.07'
$ws.Cells.Item(42, 5).Value = '  -1.65%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -4.94%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.01'
$ws.Cells.Item(44, 5).Value = '  -2.59%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '155.54'
$ws.Cells.Item(46, 5).Value = '  -1.19%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '179.59'
$ws.Cells.Item(47, 5).Value = '  -5.91%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '44.72'
$ws.Cells.Item(48, 5).Value = '  -2.68%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.132'
$ws.Cells.Item(49, 5).Value = '  +16.59%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -5.52%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '4.19'
$ws.Cells.Item(51, 5).Value = '  -2.80%  '
